# Chiffres COVID-19 Valais - update to 20.07.2020 data (adds rows for
# 44030/44031/44032, i.e. 2020-07-30/31 and 2020-08-01) and corrects a
# couple of earlier daily counts (C132, C140).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (merged A1:M1, shared string) ---------------------------------
$ws.Range("A1").Value = "Données COVID-19 Valais 20.07.2020"

# --- Corrections to existing rows (ripple through the running totals in
#     column B via the existing shared formulas) --------------------------
$ws.Range("C132").Value = 3
$ws.Range("C140").Value = 5

# --- Row 146 must inherit the "last row" outside-border formatting that
#     row 143 currently has, *before* we touch row 143's own formatting. --
$ws.Range("A143:M143").Copy()
$ws.Range("A146:M146").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 143 stops being the last row, so it switches to the regular
#     interior-row formatting that row 142 (and all rows above it) use. ---
$ws.Range("A142:M142").Copy()
$ws.Range("A143:M143").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 144 and 145 are brand new interior rows - same regular formatting.
$ws.Range("A142:M142").Copy()
$ws.Range("A144:M145").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New dates -------------------------------------------------------------
$ws.Range("A143").Value = 44029
$ws.Range("A144").Value = 44030
$ws.Range("A145").Value = 44031
$ws.Range("A146").Value = 44032

# --- Raw daily inputs (columns C, D, E, F, G, I, L, M) ---------------------
# Row 143 (2020-07-29)
$ws.Range("C143").Value = 4
$ws.Range("D143").Value = 2
$ws.Range("E143").Value = 2
$ws.Range("F143").Value = 1
$ws.Range("G143").Value = 6
$ws.Range("I143").Value = 1
$ws.Range("L143").Value = 0
$ws.Range("M143").Value = 0

# Row 144 (2020-07-30)
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 2
$ws.Range("F144").Value = 1
$ws.Range("G144").Value = 5
$ws.Range("I144").Value = 1
$ws.Range("L144").Value = 0
$ws.Range("M144").Value = 0

# Row 145 (2020-07-31)
$ws.Range("C145").Value = 1
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 2
$ws.Range("F145").Value = 1
$ws.Range("G145").Value = 5
$ws.Range("I145").Value = 0
$ws.Range("L145").Value = 0
$ws.Range("M145").Value = 0

# Row 146 (2020-08-01)
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 0
$ws.Range("E146").Value = 2
$ws.Range("F146").Value = 1
$ws.Range("G146").Value = 5
$ws.Range("I146").Value = 0
$ws.Range("L146").Value = 0
$ws.Range("M146").Value = 0

# --- Computed columns (B, H, J, K) follow the same running-total formulas
#     used throughout the sheet -------------------------------------------
$ws.Range("B143").Formula = "=B142+C143"
$ws.Range("B144").Formula = "=B143+C144"
$ws.Range("B145").Formula = "=B144+C145"
$ws.Range("B146").Formula = "=B145+C146"

$ws.Range("H143").Formula = "=G143+E143"
$ws.Range("H144").Formula = "=G144+E144"
$ws.Range("H145").Formula = "=G145+E145"
$ws.Range("H146").Formula = "=G146+E146"

$ws.Range("J143").Formula = "=J142+K143"
$ws.Range("J144").Formula = "=J143+K144"
$ws.Range("J145").Formula = "=J144+K145"
$ws.Range("J146").Formula = "=J145+K146"

$ws.Range("K143").Formula = "=L143+M143"
$ws.Range("K144").Formula = "=L144+M144"
$ws.Range("K145").Formula = "=L145+M145"
$ws.Range("K146").Formula = "=L146+M146"

# --- View bookkeeping to match the saved state of the workbook ------------
$ws.Range("A148").Select()
$ws.Application.ActiveWindow.ScrollRow = 130
